# chore: update Sheets via scheduled runner
# Refresh the computed market-price / profit columns (H, I, J, K, L, M, N)
# for the affected Leve rows across the per-job sheets, matching the
# latest pricing-tool pull. ARM!43 and CUL!137 also lose their stale
# LeveProfitHQ (N) cell entirely since there's no HQ price data anymore.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 50528.43
$ws.Range("I28").Value = 63494.125
$ws.Range("K28").Value = 63494.125
$ws.Range("M28").Value = -63009.125
# Row 33
$ws.Range("H33").Value = 48033.875
$ws.Range("I33").Value = 63790.668
$ws.Range("K33").Value = 63790.668
$ws.Range("M33").Value = -63561.668
# Row 41
$ws.Range("H41").Value = 7649.625
$ws.Range("I41").Value = 199
$ws.Range("J41").Value = 10133.167
$ws.Range("K41").Value = 199
$ws.Range("L41").Value = 10133.167
$ws.Range("M41").Value = 241
$ws.Range("N41").Value = -11013.167
# Row 53
$ws.Range("H53").Value = 2890.6365
$ws.Range("I53").Value = 76.90909000000001
$ws.Range("J53").Value = 5704.364
$ws.Range("K53").Value = 76.90909000000001
$ws.Range("L53").Value = 5704.364
$ws.Range("M53").Value = 560.09091
$ws.Range("N53").Value = -6978.364
# Row 64
$ws.Range("H64").Value = 4146.6
$ws.Range("I64").Value = 3600
$ws.Range("K64").Value = 3600
$ws.Range("M64").Value = -3352
# Row 67
$ws.Range("H67").Value = 4146.6
$ws.Range("I67").Value = 3600
$ws.Range("K67").Value = 3600
$ws.Range("M67").Value = -2742
# Row 70
$ws.Range("H70").Value = 783137.7
$ws.Range("I70").Value = 2594592.2
$ws.Range("J70").Value = 6800
$ws.Range("K70").Value = 7783776.600000001
$ws.Range("L70").Value = 20400
$ws.Range("M70").Value = -7783506.600000001
$ws.Range("N70").Value = -20940
# Row 73
$ws.Range("H73").Value = 783137.7
$ws.Range("I73").Value = 2594592.2
$ws.Range("J73").Value = 6800
$ws.Range("K73").Value = 7783776.600000001
$ws.Range("L73").Value = 20400
$ws.Range("M73").Value = -7782840.600000001
$ws.Range("N73").Value = -22272
# Row 86
$ws.Range("H86").Value = 44226.46
$ws.Range("J86").Value = 50267.637
$ws.Range("L86").Value = 50267.637
$ws.Range("N86").Value = -52513.637
# Row 89
$ws.Range("H89").Value = 44226.46
$ws.Range("J89").Value = 50267.637
$ws.Range("L89").Value = 251338.185
$ws.Range("N89").Value = -262570.185
# Row 125
$ws.Range("H125").Value = 1457.3334
$ws.Range("I125").Value = 1138.8
$ws.Range("K125").Value = 10249.2
$ws.Range("M125").Value = -7789.199999999999
# Row 138
$ws.Range("H138").Value = 3060.8667
$ws.Range("I138").Value = 2604.3635
$ws.Range("J138").Value = 3325.158
$ws.Range("K138").Value = 7813.0905
$ws.Range("L138").Value = 9975.474
$ws.Range("M138").Value = -2673.0905
$ws.Range("N138").Value = -20255.474

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 360.85
$ws.Range("I32").Value = 335.33685
$ws.Range("K32").Value = 335.33685
$ws.Range("M32").Value = -48.33685000000003
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 74
$ws.Range("H74").Value = 2170.2173
$ws.Range("I74").Value = 2021.8334
$ws.Range("K74").Value = 2021.8334
$ws.Range("M74").Value = -1147.8334
# Row 77
$ws.Range("H77").Value = 2170.2173
$ws.Range("I77").Value = 2021.8334
$ws.Range("K77").Value = 10109.167
$ws.Range("M77").Value = -5741.166999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 27105.5
$ws.Range("I99").Value = 1619.3
$ws.Range("K99").Value = 1619.3
$ws.Range("M99").Value = -121.3
# Row 105
$ws.Range("H105").Value = 2175.8333
$ws.Range("I105").Value = 2145.5557
$ws.Range("K105").Value = 2145.5557
$ws.Range("M105").Value = -398.5556999999999
# Row 107
$ws.Range("H107").Value = 8935840
$ws.Range("I107").Value = 4193.7617
$ws.Range("K107").Value = 4193.7617
$ws.Range("M107").Value = -2273.7617

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1900.238
$ws.Range("I58").Value = 1440.5714
$ws.Range("J58").Value = 2130.0715
$ws.Range("K58").Value = 1440.5714
$ws.Range("L58").Value = 2130.0715
$ws.Range("M58").Value = -1237.5714
$ws.Range("N58").Value = -2536.0715
# Row 104
$ws.Range("H104").Value = 59888
$ws.Range("J104").Value = 59888
$ws.Range("L104").Value = 59888
$ws.Range("N104").Value = -65130
# Row 132
$ws.Range("H132").Value = 2728.3809
$ws.Range("I132").Value = 1825.2222
$ws.Range("K132").Value = 5475.6666
$ws.Range("M132").Value = -2945.6666
# Row 134
$ws.Range("H134").Value = 2478.0967
$ws.Range("I134").Value = 1748.8422
$ws.Range("K134").Value = 5246.5266
$ws.Range("M134").Value = -2711.5266
# Row 136
$ws.Range("H136").Value = 1900.238
$ws.Range("I136").Value = 1440.5714
$ws.Range("J136").Value = 2130.0715
$ws.Range("K136").Value = 4321.7142
$ws.Range("L136").Value = 6390.2145
$ws.Range("M136").Value = -1771.7142
$ws.Range("N136").Value = -11490.2145

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 4011.75
$ws.Range("I8").Value = 4011.75
$ws.Range("K8").Value = 12035.25
$ws.Range("M8").Value = -11896.25
# Row 23
$ws.Range("H23").Value = 1310.7273
$ws.Range("I23").Value = 775
$ws.Range("K23").Value = 2325
$ws.Range("M23").Value = -2090
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 140
$ws.Range("H140").Value = 1814
$ws.Range("I140").Value = 1244.25
$ws.Range("K140").Value = 3732.75
$ws.Range("M140").Value = 1447.25

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 8061.3335
$ws.Range("I126").Value = 10531.167
$ws.Range("J126").Value = 4768.222
$ws.Range("K126").Value = 31593.501
$ws.Range("L126").Value = 14304.666
$ws.Range("M126").Value = -29123.501
$ws.Range("N126").Value = -19244.666
# Row 132
$ws.Range("H132").Value = 316099.53
$ws.Range("I132").Value = 419801.5
$ws.Range("J132").Value = 4993.625
$ws.Range("K132").Value = 1259404.5
$ws.Range("L132").Value = 14980.875
$ws.Range("M132").Value = -1256874.5
$ws.Range("N132").Value = -20040.875

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3565.5334
$ws.Range("I40").Value = 3043.9092
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3043.9092
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2907.9092
$ws.Range("N40").Value = -5272
# Row 93
$ws.Range("H93").Value = 1036.619
$ws.Range("I93").Value = 933.2778
$ws.Range("K93").Value = 933.2778
$ws.Range("M93").Value = 314.7222
# Row 136
$ws.Range("H136").Value = 9782.583000000001
$ws.Range("I136").Value = 4450.75
$ws.Range("J136").Value = 12448.5
$ws.Range("K136").Value = 13352.25
$ws.Range("L136").Value = 37345.5
$ws.Range("M136").Value = -10802.25
$ws.Range("N136").Value = -42445.5

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 387.36365
$ws.Range("I107").Value = 384.35715
$ws.Range("J107").Value = 392.625
$ws.Range("K107").Value = 1153.07145
$ws.Range("L107").Value = 1177.875
$ws.Range("M107").Value = 766.9285500000001
$ws.Range("N107").Value = -5017.875
# Row 132
$ws.Range("H132").Value = 259409.16
$ws.Range("I132").Value = 305892.7
$ws.Range("K132").Value = 917678.1000000001
$ws.Range("M132").Value = -915148.1000000001
